$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 45

# Build the new date label in a scratch cell formatted as text (so Excel
# does not auto-convert the "01-08-2021" string into a date serial), then
# copy just the resulting value into the target cell and clean up the
# scratch cell again.
$scratch = $ws.Cells.Item(100, 26)
$scratch.NumberFormat = "@"
$scratch.Value = "01-08-2021"
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($row, 2).Value = 110.15
$ws.Cells.Item($row, 3).Value = 108.41
$ws.Cells.Item($row, 4).Value = 111.67
$ws.Cells.Item($row, 5).Value = 108.36
$ws.Cells.Item($row, 6).Value = 117.88
